# Insert a new weekly record at row 241 (Vega Modelo de Temuco - Albahaca),
# pushing the existing rows 241:269 down to 242:270.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row 241 (copies formatting
# from the row above, same as Excel's default "Insert" behaviour for a
# right-click > Insert on a selected row).
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row with the new weekly observation. It is
# identical to the record that used to be on row 241 (now row 242) except
# for a newer reporting date (Fecha).
$ws.Range("A241").Value = 10
$ws.Range("B241").Value = "Vega Modelo de Temuco"
$ws.Range("C241").Value = "La Araucanía"
$ws.Range("D241").Value = 44816
$ws.Range("E241").Value = 9
$ws.Range("F241").Value = 100112052
$ws.Range("G241").Value = "Albahaca"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 50
$ws.Range("K241").Value = 6000
$ws.Range("L241").Value = 6000
$ws.Range("M241").Value = 6000
$ws.Range("N241").Value = "$/paquete"
$ws.Range("O241").Value = "Región de Arica y Parinacota"
$ws.Range("P241").Value = 6000
$ws.Range("Q241").Value = 1
$ws.Range("R241").Value = "Hortaliza"
